$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-populated Train Final Accuracy (D) and Test Final
# Accuracy (F) values for the "2nd" change-id rows (4 and 5).
$ws.Range("D4").Value = 67.64
$ws.Range("F4").Value = 54.88
$ws.Range("D5").Value = 67.64
$ws.Range("F5").Value = 54.88

# Zoom out the sheet view and move the selection to G5, matching the
# updated sheetView/selection state.
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("G5").Select()
